$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Field" names to lowercase variants (fr_per_current_mean / fr_per_current_sem)
$ws.Range("A7").Value = "fr_per_current_mean"
$ws.Range("A8").Value = "fr_per_current_sem"

# Update the active cell selection to C7
$null = $ws.Range("C7").Select()
